# AccountCreationTestData.xlsx update:
#  - rename sheet "AccountCreationData" -> "accountCreationTest"
#  - update the three test email addresses on that sheet
#  - make that sheet the active / selected one (tab + A4 selection),
#    moving activation away from "verifyTotalPriceTest"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AccountCreationData")
$ws.Name = "accountCreationTest"

$ws.Range("A2").Value = "nht2@gmail.com"
$ws.Range("A3").Value = "qhsd21@gmail.com"
$ws.Range("A4").Value = "jt1@gmail.com"

$ws.Activate()
$ws.Range("A4").Select()
